# prépa - traitement données erié élèves
# mise en forme du database et du rmd pour générer les graphs emo et cps (experimental)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header columns G1:L1 (new variables tracked for the study) ---
$ws.Range("G1").Value = "code_temps_1"
$ws.Range("H1").Value = "code_temps_2"
$ws.Range("I1").Value = "signalé le"
$ws.Range("J1").Value = "par"
$ws.Range("K1").Value = "prétest envoyé le"
$ws.Range("L1").Value = "post-test envoyé le"

# --- The old single sample/example row is wiped (only C2 keeps its style) ---
$ws.Range("A2:F2").ClearContents()

# --- Widen the new date columns K:L ---
$ws.Range("K1:L1").ColumnWidth = 21.166666666666668

# --- Selection ends up on the newly added K1:L1 header cells ---
$ws.Range("K1:L1").Select()

# --- Add the explanatory text box (legend for the column codes) ---
$shp = $ws.Shapes.AddTextbox(1, 449, 31, 296.857165, 254.428583)
$shp.Name = "ZoneTexte 1"

$legend = "classe (identifiant de la classe) : a b c d e ... aa ab .. zz`n" +
  "gest (titulariat) : solo ou duo`n" +
  "deg (degré suivis) : 1-2 ... 8`n" +
  "cyc (cycle) : 1 ou 2`n" +
  "sit (direction) : Haut-Lac, Collombey, Monthey, Fully, Saxon`n" +
  "code_temps_1 : code à entrer par profs et élèves au temps 1`n" +
  "code_temps_2 : code à entrer par profs et élèves au temps 2`n" +
  "signalé le : date d'inscription`n" +
  "par : erie concerné·e`n" +
  "`n" +
  "Construction du code :`n" +
  "- commence par 123`n" +
  "- se poursuit par code classe`n" +
  "- se poursuit par temps de mesure 1 ou 2`n" +
  "- finit par 5 `n" +
  "`n" +
  "Exemple pour classe a au temps 1 : 123a15"

$shp.TextFrame2.TextRange.Text = $legend
$shp.TextFrame2.TextRange.Font.Size = 11
$shp.TextFrame2.TextRange.Font.Name = "Calibri"
